$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.062008333333333
$ws.Range("H2").Value = 9.186025
$ws.Range("I2").Value = 0.762008591445137
$ws.Range("J2").Value = 0.7620085914451371
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.22092
$ws.Range("N2").Value = 0.66276
$ws.Range("O2").Value = 0.06768692722940421
$ws.Range("P2").Value = 0.06768692722940423
$ws.Range("Q2").Value = 0.676458881
$ws.Range("R2").Value = 6.088129929000001
$ws.Range("S2").Value = 0.05157802007732779
$ws.Range("T2").Value = 0.05157802007732781
$ws.Range("G3").Value = 3.062008333333333
$ws.Range("H3").Value = 9.186025
$ws.Range("I3").Value = 0.762008591445137
$ws.Range("J3").Value = 0.7620085914451371
$ws.Range("O3").Value = 0.193804638305004
$ws.Range("P3").Value = 0.193804638305004
$ws.Range("Q3").Value = 1.9368713299111109
$ws.Range("R3").Value = 17.4318419692
$ws.Range("S3").Value = 0.14768079945033033
$ws.Range("T3").Value = 0.14768079945033036
$ws.Range("G4").Value = 3.062008333333333
$ws.Range("H4").Value = 9.186025
$ws.Range("I4").Value = 0.762008591445137
$ws.Range("J4").Value = 0.7620085914451371
$ws.Range("M4").Value = 0.1410223333333333
$ws.Range("N4").Value = 0.423067
$ws.Range("O4").Value = 0.04320735295153955
$ws.Range("P4").Value = 0.04320735295153956
$ws.Range("Q4").Value = 0.43181155985277764
$ws.Range("R4").Value = 3.8863040386750005
$ws.Range("S4").Value = 0.03292437416267554
$ws.Range("T4").Value = 0.032924374162675546
$ws.Range("G5").Value = 3.062008333333333
$ws.Range("H5").Value = 9.186025
$ws.Range("I5").Value = 0.762008591445137
$ws.Range("J5").Value = 0.7620085914451371
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.114977
$ws.Range("N5").Value = 0.344931
$ws.Range("O5").Value = 0.03522741187785264
$ws.Range("P5").Value = 0.03522741187785266
$ws.Range("Q5").Value = 0.3520605321416666
$ws.Range("R5").Value = 3.1685447892750003
$ws.Range("S5").Value = 0.026843590505300177
$ws.Range("T5").Value = 0.026843590505300198
$ws.Range("G6").Value = 3.062008333333333
$ws.Range("H6").Value = 9.186025
$ws.Range("I6").Value = 0.762008591445137
$ws.Range("J6").Value = 0.7620085914451371
$ws.Range("M6").Value = 2.032918666666667
$ws.Range("N6").Value = 6.098756
$ws.Range("O6").Value = 0.6228590342837411
$ws.Range("P6").Value = 0.6228590342837412
$ws.Range("Q6").Value = 6.224813898322223
$ws.Range("R6").Value = 56.0233250849
$ws.Range("S6").Value = 0.47462393538343184
$ws.Range("T6").Value = 0.474623935383432
$ws.Range("G7").Value = 3.062008333333333
$ws.Range("H7").Value = 9.186025
$ws.Range("I7").Value = 0.762008591445137
$ws.Range("J7").Value = 0.7620085914451371
$ws.Range("M7").Value = 0.121463
$ws.Range("N7").Value = 0.364389
$ws.Range("O7").Value = 0.03721463535245846
$ws.Range("P7").Value = 0.03721463535245847
$ws.Range("Q7").Value = 0.3719207181916666
$ws.Range("R7").Value = 3.3472864637250006
$ws.Range("S7").Value = 0.02835787186607127
$ws.Range("T7").Value = 0.02835787186607128
$ws.Range("I8").Value = 0.007987298232312442
$ws.Range("J8").Value = 0.007987298232312444
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 0.22092
$ws.Range("N8").Value = 0.66276
$ws.Range("O8").Value = 0.06768692722940421
$ws.Range("P8").Value = 0.06768692722940423
$ws.Range("Q8").Value = 0.00709057468
$ws.Range("R8").Value = 0.06381517212000001
$ws.Range("S8").Value = 0.0005406356742100811
$ws.Range("T8").Value = 0.0005406356742100813
$ws.Range("I9").Value = 0.007987298232312442
$ws.Range("J9").Value = 0.007987298232312444
$ws.Range("O9").Value = 0.193804638305004
$ws.Range("P9").Value = 0.193804638305004
$ws.Range("S9").Value = 0.0015479754449475107
$ws.Range("T9").Value = 0.001547975444947511
$ws.Range("I10").Value = 0.007987298232312442
$ws.Range("J10").Value = 0.007987298232312444
$ws.Range("M10").Value = 0.1410223333333333
$ws.Range("N10").Value = 0.423067
$ws.Range("O10").Value = 0.04320735295153955
$ws.Range("P10").Value = 0.04320735295153956
$ws.Range("Q10").Value = 0.0045262058032222215
$ws.Range("R10").Value = 0.04073585222900001
$ws.Range("S10").Value = 0.00034511001385273166
$ws.Range("T10").Value = 0.00034511001385273177
$ws.Range("I11").Value = 0.007987298232312442
$ws.Range("J11").Value = 0.007987298232312444
$ws.Range("K11").Value = 2.0
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.114977
$ws.Range("N11").Value = 0.344931
$ws.Range("O11").Value = 0.03522741187785264
$ws.Range("P11").Value = 0.03522741187785266
$ws.Range("Q11").Value = 0.0036902634663333335
$ws.Range("R11").Value = 0.033212371197000004
$ws.Range("S11").Value = 0.0002813718446209147
$ws.Range("T11").Value = 0.0002813718446209149
$ws.Range("I12").Value = 0.007987298232312442
$ws.Range("J12").Value = 0.007987298232312444
$ws.Range("M12").Value = 2.032918666666667
$ws.Range("N12").Value = 6.098756
$ws.Range("O12").Value = 0.6228590342837411
$ws.Range("P12").Value = 0.6228590342837412
$ws.Range("Q12").Value = 0.0652478798857778
$ws.Range("R12").Value = 0.5872309189720001
$ws.Range("S12").Value = 0.00497496086351436
$ws.Range("T12").Value = 0.004974960863514362
$ws.Range("I13").Value = 0.007987298232312442
$ws.Range("J13").Value = 0.007987298232312444
$ws.Range("M13").Value = 0.121463
$ws.Range("N13").Value = 0.364389
$ws.Range("O13").Value = 0.03721463535245846
$ws.Range("P13").Value = 0.03721463535245847
$ws.Range("Q13").Value = 0.0038984359603333337
$ws.Range("R13").Value = 0.035085923643000005
$ws.Range("S13").Value = 0.0002972443911668436
$ws.Range("T13").Value = 0.0002972443911668437
$ws.Range("G14").Value = 0.02852133333333333
$ws.Range("H14").Value = 0.085564
$ws.Range("I14").Value = 0.007097792910253532
$ws.Range("J14").Value = 0.007097792910253533
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 0.22092
$ws.Range("N14").Value = 0.66276
$ws.Range("O14").Value = 0.06768692722940421
$ws.Range("P14").Value = 0.06768692722940423
$ws.Range("Q14").Value = 0.006300932959999999
$ws.Range("R14").Value = 0.056708396640000006
$ws.Range("S14").Value = 0.00048042779220571196
$ws.Range("T14").Value = 0.0004804277922057121
$ws.Range("G15").Value = 0.02852133333333333
$ws.Range("H15").Value = 0.085564
$ws.Range("I15").Value = 0.007097792910253532
$ws.Range("J15").Value = 0.007097792910253533
$ws.Range("O15").Value = 0.193804638305004
$ws.Range("P15").Value = 0.193804638305004
$ws.Range("Q15").Value = 0.018041150385777773
$ws.Range("R15").Value = 0.162370353472
$ws.Range("S15").Value = 0.0013755851877355074
$ws.Range("T15").Value = 0.0013755851877355076
$ws.Range("G16").Value = 0.02852133333333333
$ws.Range("H16").Value = 0.085564
$ws.Range("I16").Value = 0.007097792910253532
$ws.Range("J16").Value = 0.007097792910253533
$ws.Range("M16").Value = 0.1410223333333333
$ws.Range("N16").Value = 0.423067
$ws.Range("O16").Value = 0.04320735295153955
$ws.Range("P16").Value = 0.04320735295153956
$ws.Range("Q16").Value = 0.004022144976444443
$ws.Range("R16").Value = 0.036199304788
$ws.Range("S16").Value = 0.0003066768434502594
$ws.Range("T16").Value = 0.00030667684345025954
$ws.Range("G17").Value = 0.02852133333333333
$ws.Range("H17").Value = 0.085564
$ws.Range("I17").Value = 0.007097792910253532
$ws.Range("J17").Value = 0.007097792910253533
$ws.Range("K17").Value = 2.0
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.114977
$ws.Range("N17").Value = 0.344931
$ws.Range("O17").Value = 0.03522741187785264
$ws.Range("P17").Value = 0.03522741187785266
$ws.Range("Q17").Value = 0.003279297342666666
$ws.Range("R17").Value = 0.029513676084
$ws.Range("S17").Value = 0.0002500368742732035
$ws.Range("T17").Value = 0.00025003687427320367
$ws.Range("G18").Value = 0.02852133333333333
$ws.Range("H18").Value = 0.085564
$ws.Range("I18").Value = 0.007097792910253532
$ws.Range("J18").Value = 0.007097792910253533
$ws.Range("M18").Value = 2.032918666666667
$ws.Range("N18").Value = 6.098756
$ws.Range("O18").Value = 0.6228590342837411
$ws.Range("P18").Value = 0.6228590342837412
$ws.Range("Q18").Value = 0.05798155093155556
$ws.Range("R18").Value = 0.521833958384
$ws.Range("S18").Value = 0.004420924437626499
$ws.Range("T18").Value = 0.004420924437626501
$ws.Range("G19").Value = 0.02852133333333333
$ws.Range("H19").Value = 0.085564
$ws.Range("I19").Value = 0.007097792910253532
$ws.Range("J19").Value = 0.007097792910253533
$ws.Range("M19").Value = 0.121463
$ws.Range("N19").Value = 0.364389
$ws.Range("O19").Value = 0.03721463535245846
$ws.Range("P19").Value = 0.03721463535245847
$ws.Range("Q19").Value = 0.003464286710666666
$ws.Range("R19").Value = 0.031178580396
$ws.Range("S19").Value = 0.0002641417749623501
$ws.Range("T19").Value = 0.00026414177496235017
$ws.Range("G20").Value = 0.1165276666666667
$ws.Range("H20").Value = 0.349583
$ws.Range("I20").Value = 0.02899896847909355
$ws.Range("J20").Value = 0.02899896847909356
$ws.Range("K20").Value = 3.0
$ws.Range("L20").Value = 1.0
$ws.Range("M20").Value = 0.22092
$ws.Range("N20").Value = 0.66276
$ws.Range("O20").Value = 0.06768692722940421
$ws.Range("P20").Value = 0.06768692722940423
$ws.Range("Q20").Value = 0.025743292120000007
$ws.Range("R20").Value = 0.23168962908
$ws.Range("S20").Value = 0.0019628510691721917
$ws.Range("T20").Value = 0.0019628510691721925
$ws.Range("G21").Value = 0.1165276666666667
$ws.Range("H21").Value = 0.349583
$ws.Range("I21").Value = 0.02899896847909355
$ws.Range("J21").Value = 0.02899896847909356
$ws.Range("O21").Value = 0.193804638305004
$ws.Range("P21").Value = 0.193804638305004
$ws.Range("Q21").Value = 0.07370949786488891
$ws.Range("R21").Value = 0.663385480784
$ws.Range("S21").Value = 0.005620134597308937
$ws.Range("T21").Value = 0.005620134597308939
$ws.Range("G22").Value = 0.1165276666666667
$ws.Range("H22").Value = 0.349583
$ws.Range("I22").Value = 0.02899896847909355
$ws.Range("J22").Value = 0.02899896847909356
$ws.Range("M22").Value = 0.1410223333333333
$ws.Range("N22").Value = 0.423067
$ws.Range("O22").Value = 0.04320735295153955
$ws.Range("P22").Value = 0.04320735295153956
$ws.Range("Q22").Value = 0.016433003451222224
$ws.Range("R22").Value = 0.147897031061
$ws.Range("S22").Value = 0.001252968666306765
$ws.Range("T22").Value = 0.0012529686663067658
$ws.Range("G23").Value = 0.1165276666666667
$ws.Range("H23").Value = 0.349583
$ws.Range("I23").Value = 0.02899896847909355
$ws.Range("J23").Value = 0.02899896847909356
$ws.Range("K23").Value = 2.0
$ws.Range("L23").Value = 0.6666666666666666
$ws.Range("M23").Value = 0.114977
$ws.Range("N23").Value = 0.344931
$ws.Range("O23").Value = 0.03522741187785264
$ws.Range("P23").Value = 0.03522741187785266
$ws.Range("Q23").Value = 0.013398001530333336
$ws.Range("R23").Value = 0.12058201377299999
$ws.Range("S23").Value = 0.0010215586066458944
$ws.Range("T23").Value = 0.0010215586066458953
$ws.Range("G24").Value = 0.1165276666666667
$ws.Range("H24").Value = 0.349583
$ws.Range("I24").Value = 0.02899896847909355
$ws.Range("J24").Value = 0.02899896847909356
$ws.Range("M24").Value = 2.032918666666667
$ws.Range("N24").Value = 6.098756
$ws.Range("O24").Value = 0.6228590342837411
$ws.Range("P24").Value = 0.6228590342837412
$ws.Range("Q24").Value = 0.2368912687497779
$ws.Range("R24").Value = 2.1320214187479998
$ws.Range("S24").Value = 0.018062269502112857
$ws.Range("T24").Value = 0.018062269502112867
$ws.Range("G25").Value = 0.1165276666666667
$ws.Range("H25").Value = 0.349583
$ws.Range("I25").Value = 0.02899896847909355
$ws.Range("J25").Value = 0.02899896847909356
$ws.Range("M25").Value = 0.121463
$ws.Range("N25").Value = 0.364389
$ws.Range("O25").Value = 0.03721463535245846
$ws.Range("P25").Value = 0.03721463535245847
$ws.Range("Q25").Value = 0.014153799976333337
$ws.Range("R25").Value = 0.127384199787
$ws.Range("S25").Value = 0.0010791860375469033
$ws.Range("T25").Value = 0.001079186037546904
$ws.Range("G26").Value = 0.7333496666666667
$ws.Range("H26").Value = 2.200049
$ws.Range("I26").Value = 0.1825007268759101
$ws.Range("J26").Value = 0.1825007268759101
$ws.Range("K26").Value = 3.0
$ws.Range("L26").Value = 1.0
$ws.Range("M26").Value = 0.22092
$ws.Range("N26").Value = 0.66276
$ws.Range("O26").Value = 0.06768692722940421
$ws.Range("P26").Value = 0.06768692722940423
$ws.Range("Q26").Value = 0.16201160836
$ws.Range("R26").Value = 1.45810447524
$ws.Range("S26").Value = 0.0123529134193631
$ws.Range("T26").Value = 0.012352913419363103
$ws.Range("G27").Value = 0.7333496666666667
$ws.Range("H27").Value = 2.200049
$ws.Range("I27").Value = 0.1825007268759101
$ws.Range("J27").Value = 0.1825007268759101
$ws.Range("O27").Value = 0.193804638305004
$ws.Range("P27").Value = 0.193804638305004
$ws.Range("Q27").Value = 0.4638798427502222
$ws.Range("R27").Value = 4.174918584752
$ws.Range("S27").Value = 0.03536948736258608
$ws.Range("T27").Value = 0.03536948736258608
$ws.Range("G28").Value = 0.7333496666666667
$ws.Range("H28").Value = 2.200049
$ws.Range("I28").Value = 0.1825007268759101
$ws.Range("J28").Value = 0.1825007268759101
$ws.Range("M28").Value = 0.1410223333333333
$ws.Range("N28").Value = 0.423067
$ws.Range("O28").Value = 0.04320735295153955
$ws.Range("P28").Value = 0.04320735295153956
$ws.Range("Q28").Value = 0.10341868114255554
$ws.Range("R28").Value = 0.930768130283
$ws.Range("S28").Value = 0.007885373320039968
$ws.Range("T28").Value = 0.00788537332003997
$ws.Range("G29").Value = 0.7333496666666667
$ws.Range("H29").Value = 2.200049
$ws.Range("I29").Value = 0.1825007268759101
$ws.Range("J29").Value = 0.1825007268759101
$ws.Range("K29").Value = 2.0
$ws.Range("L29").Value = 0.6666666666666666
$ws.Range("M29").Value = 0.114977
$ws.Range("N29").Value = 0.344931
$ws.Range("O29").Value = 0.03522741187785264
$ws.Range("P29").Value = 0.03522741187785266
$ws.Range("Q29").Value = 0.08431834462433334
$ws.Range("R29").Value = 0.7588651016189999
$ws.Range("S29").Value = 0.006429028273665175
$ws.Range("T29").Value = 0.006429028273665179
$ws.Range("G30").Value = 0.7333496666666667
$ws.Range("H30").Value = 2.200049
$ws.Range("I30").Value = 0.1825007268759101
$ws.Range("J30").Value = 0.1825007268759101
$ws.Range("M30").Value = 2.032918666666667
$ws.Range("N30").Value = 6.098756
$ws.Range("O30").Value = 0.6228590342837411
$ws.Range("P30").Value = 0.6228590342837412
$ws.Range("Q30").Value = 1.4908402265604448
$ws.Range("R30").Value = 13.417562039043998
$ws.Range("S30").Value = 0.11367222649801016
$ws.Range("T30").Value = 0.11367222649801018
$ws.Range("G31").Value = 0.7333496666666667
$ws.Range("H31").Value = 2.200049
$ws.Range("I31").Value = 0.1825007268759101
$ws.Range("J31").Value = 0.1825007268759101
$ws.Range("M31").Value = 0.121463
$ws.Range("N31").Value = 0.364389
$ws.Range("O31").Value = 0.03721463535245846
$ws.Range("P31").Value = 0.03721463535245847
$ws.Range("Q31").Value = 0.08907485056233333
$ws.Range("R31").Value = 0.801673655061
$ws.Range("S31").Value = 0.00679169800224561
$ws.Range("T31").Value = 0.006791698002245611
$ws.Range("E32").Value = 3.0
$ws.Range("F32").Value = 1.0
$ws.Range("G32").Value = 0.04583566666666666
$ws.Range("H32").Value = 0.137507
$ws.Range("I32").Value = 0.01140662205729316
$ws.Range("J32").Value = 0.01140662205729317
$ws.Range("K32").Value = 3.0
$ws.Range("L32").Value = 1.0
$ws.Range("M32").Value = 0.22092
$ws.Range("N32").Value = 0.66276
$ws.Range("O32").Value = 0.06768692722940421
$ws.Range("P32").Value = 0.06768692722940423
$ws.Range("Q32").Value = 0.01012601548
$ws.Range("R32").Value = 0.09113413932
$ws.Range("S32").Value = 0.0007720791971253191
$ws.Range("T32").Value = 0.0007720791971253199
$ws.Range("E33").Value = 3.0
$ws.Range("F33").Value = 1.0
$ws.Range("G33").Value = 0.04583566666666666
$ws.Range("H33").Value = 0.137507
$ws.Range("I33").Value = 0.01140662205729316
$ws.Range("J33").Value = 0.01140662205729317
$ws.Range("O33").Value = 0.193804638305004
$ws.Range("P33").Value = 0.193804638305004
$ws.Range("Q33").Value = 0.028993320392888884
$ws.Range("R33").Value = 0.260939883536
$ws.Range("S33").Value = 0.0022106562620955814
$ws.Range("T33").Value = 0.0022106562620955836
$ws.Range("E34").Value = 3.0
$ws.Range("F34").Value = 1.0
$ws.Range("G34").Value = 0.04583566666666666
$ws.Range("H34").Value = 0.137507
$ws.Range("I34").Value = 0.01140662205729316
$ws.Range("J34").Value = 0.01140662205729317
$ws.Range("M34").Value = 0.1410223333333333
$ws.Range("N34").Value = 0.423067
$ws.Range("O34").Value = 0.04320735295153955
$ws.Range("P34").Value = 0.04320735295153956
$ws.Range("Q34").Value = 0.00646385266322222
$ws.Range("R34").Value = 0.058174673969
$ws.Range("S34").Value = 0.0004928499452142818
$ws.Range("T34").Value = 0.0004928499452142823
$ws.Range("E35").Value = 3.0
$ws.Range("F35").Value = 1.0
$ws.Range("G35").Value = 0.04583566666666666
$ws.Range("H35").Value = 0.137507
$ws.Range("I35").Value = 0.01140662205729316
$ws.Range("J35").Value = 0.01140662205729317
$ws.Range("K35").Value = 2.0
$ws.Range("L35").Value = 0.6666666666666666
$ws.Range("M35").Value = 0.114977
$ws.Range("N35").Value = 0.344931
$ws.Range("O35").Value = 0.03522741187785264
$ws.Range("P35").Value = 0.03522741187785266
$ws.Range("Q35").Value = 0.005270047446333333
$ws.Range("R35").Value = 0.047430427016999996
$ws.Range("S35").Value = 0.00040182577334726494
$ws.Range("T35").Value = 0.00040182577334726553
$ws.Range("E36").Value = 3.0
$ws.Range("F36").Value = 1.0
$ws.Range("G36").Value = 0.04583566666666666
$ws.Range("H36").Value = 0.137507
$ws.Range("I36").Value = 0.01140662205729316
$ws.Range("J36").Value = 0.01140662205729317
$ws.Range("M36").Value = 2.032918666666667
$ws.Range("N36").Value = 6.098756
$ws.Range("O36").Value = 0.6228590342837411
$ws.Range("P36").Value = 0.6228590342837412
$ws.Range("Q36").Value = 0.0931801823657778
$ws.Range("R36").Value = 0.8386216412919999
$ws.Range("S36").Value = 0.007104717599045238
$ws.Range("T36").Value = 0.007104717599045246
$ws.Range("E37").Value = 3.0
$ws.Range("F37").Value = 1.0
$ws.Range("G37").Value = 0.04583566666666666
$ws.Range("H37").Value = 0.137507
$ws.Range("I37").Value = 0.01140662205729316
$ws.Range("J37").Value = 0.01140662205729317
$ws.Range("M37").Value = 0.121463
$ws.Range("N37").Value = 0.364389
$ws.Range("O37").Value = 0.03721463535245846
$ws.Range("P37").Value = 0.03721463535245847
$ws.Range("Q37").Value = 0.005567337580333333
$ws.Range("R37").Value = 0.050106038223
$ws.Range("S37").Value = 0.0004244932804654745
$ws.Range("T37").Value = 0.00042449328046547497
